# NA OCP Result_Test.xlsx - "Changes of NAM for Prod"
# Updates PickupID/POD No (column C) and Fail Log (column F) values on the
# "Result" sheet to reflect the latest automation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (PickupID/POD No) updates ---
$ws.Range("C2").Value = "10294941"
$ws.Range("C3").Value = "10294944"
$ws.Range("C4").Value = "10294812"
$ws.Range("C5").Value = "10294830"
$ws.Range("C11").Value = "10294915"
$ws.Range("C12").Value = "10294856"
$ws.Range("C13").Value = "10294889"
$ws.Range("C14").Value = "10294905"
$ws.Range("C24").Value = "137213748"

# --- Column F (Fail Log) updates ---
$f3 = 'Cannot invoke "org.openqa.selenium.WebElement.getText()" because the return value of "connect_OCBaseMethods.TCAcknowledge.isElementPresent(String)" is null'
$ws.Range("F3").Value = $f3

$f26 = @'
Expected condition failed: waiting for visibility of [[ChromeDriver: chrome on WINDOWS (678a429e6d9e620b861fd0a924b694ab)] -> id: GreyTick] (tried for 60 second(s) with 500 milliseconds interval)
Build info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'
System info: host: 'SIPL92', ip: '10.212.130.60', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 122.0.6261.113, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:53284}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 678a429e6d9e620b861fd0a924b694ab
'@
$ws.Range("F26").Value = $f26
